$d = $word.ActiveDocument
$d.Content.Find.Execute("2023-07-06 Thursday", $true, $true, $false, $false, $false, $true, 1, $false, "2023-07-07 Friday", 2) | Out-Null
$d.Content.Find.Execute("99-12=", $true, $true, $false, $false, $false, $true, 1, $false, "29-11=", 2) | Out-Null
$d.Content.Find.Execute("17-10=", $true, $true, $false, $false, $false, $true, 1, $false, "96-24=", 2) | Out-Null
$d.Content.Find.Execute("68-16=", $true, $true, $false, $false, $false, $true, 1, $false, "46-24=", 2) | Out-Null
$d.Content.Find.Execute("35+7=", $true, $true, $false, $false, $false, $true, 1, $false, "15+33=", 2) | Out-Null
$d.Content.Find.Execute("27+31=", $true, $true, $false, $false, $false, $true, 1, $false, "38-29=", 2) | Out-Null
$d.Content.Find.Execute("66-41=", $true, $true, $false, $false, $false, $true, 1, $false, "17+19=", 2) | Out-Null
$d.Content.Find.Execute("94-60=", $true, $true, $false, $false, $false, $true, 1, $false, "35+46=", 2) | Out-Null
$d.Content.Find.Execute("61-4=", $true, $true, $false, $false, $false, $true, 1, $false, "31-6=", 2) | Out-Null
$d.Content.Find.Execute("64+20=", $true, $true, $false, $false, $false, $true, 1, $false, "34-16=", 2) | Out-Null
$d.Content.Find.Execute("45+16=", $true, $true, $false, $false, $false, $true, 1, $false, "8+2=", 2) | Out-Null
$d.Content.Find.Execute("8-4=", $true, $true, $false, $false, $false, $true, 1, $false, "5+76=", 2) | Out-Null
$d.Content.Find.Execute("40+5=", $true, $true, $false, $false, $false, $true, 1, $false, "26+6=", 2) | Out-Null
$d.Content.Find.Execute("68-56=", $true, $true, $false, $false, $false, $true, 1, $false, "50+22=", 2) | Out-Null
$d.Content.Find.Execute("66-11=", $true, $true, $false, $false, $false, $true, 1, $false, "96-78=", 2) | Out-Null
$d.Content.Find.Execute("50-2=", $true, $true, $false, $false, $false, $true, 1, $false, "8+28=", 2) | Out-Null
$d.Content.Find.Execute("9+75=", $true, $true, $false, $false, $false, $true, 1, $false, "36+62=", 2) | Out-Null
$d.Content.Find.Execute("73+5=", $true, $true, $false, $false, $false, $true, 1, $false, "7+89=", 2) | Out-Null
$d.Content.Find.Execute("14+45=", $true, $true, $false, $false, $false, $true, 1, $false, "25+37=", 2) | Out-Null
$d.Content.Find.Execute("68-46=", $true, $true, $false, $false, $false, $true, 1, $false, "41-33=", 2) | Out-Null
$d.Content.Find.Execute("82-30=", $true, $true, $false, $false, $false, $true, 1, $false, "65-9=", 2) | Out-Null
$d.Content.Find.Execute("34-15=", $true, $true, $false, $false, $false, $true, 1, $false, "18+36=", 2) | Out-Null
$d.Content.Find.Execute("79-31=", $true, $true, $false, $false, $false, $true, 1, $false, "67+12=", 2) | Out-Null
$d.Content.Find.Execute("9+21=", $true, $true, $false, $false, $false, $true, 1, $false, "4+9=", 2) | Out-Null
$d.Content.Find.Execute("6+41=", $true, $true, $false, $false, $false, $true, 1, $false, "1+72=", 2) | Out-Null
$d.Content.Find.Execute("49-22=", $true, $true, $false, $false, $false, $true, 1, $false, "73-51=", 2) | Out-Null
$d.Content.Find.Execute("56+5=", $true, $true, $false, $false, $false, $true, 1, $false, "20+23=", 2) | Out-Null
$d.Content.Find.Execute("55-34=", $true, $true, $false, $false, $false, $true, 1, $false, "59+14=", 2) | Out-Null
$d.Content.Find.Execute("86-39=", $true, $true, $false, $false, $false, $true, 1, $false, "92-77=", 2) | Out-Null
$d.Content.Find.Execute("41+7=", $true, $true, $false, $false, $false, $true, 1, $false, "76+1=", 2) | Out-Null
$d.Content.Find.Execute("62-30=", $true, $true, $false, $false, $false, $true, 1, $false, "48-14=", 2) | Out-Null
$d.Content.Find.Execute("77-24=", $true, $true, $false, $false, $false, $true, 1, $false, "17+54=", 2) | Out-Null
$d.Content.Find.Execute("7+81=", $true, $true, $false, $false, $false, $true, 1, $false, "52+47=", 2) | Out-Null
$d.Content.Find.Execute("44+14=", $true, $true, $false, $false, $false, $true, 1, $false, "48-38=", 2) | Out-Null
$d.Content.Find.Execute("45-43=", $true, $true, $false, $false, $false, $true, 1, $false, "19+50=", 2) | Out-Null
$d.Content.Find.Execute("92+2=", $true, $true, $false, $false, $false, $true, 1, $false, "81-70=", 2) | Out-Null
$d.Content.Find.Execute("90-24=", $true, $true, $false, $false, $false, $true, 1, $false, "75-61=", 2) | Out-Null
$d.Content.Find.Execute("39+37=", $true, $true, $false, $false, $false, $true, 1, $false, "35+49=", 2) | Out-Null
$d.Content.Find.Execute("43-22=", $true, $true, $false, $false, $false, $true, 1, $false, "81-59=", 2) | Out-Null
$d.Content.Find.Execute("1+20=", $true, $true, $false, $false, $false, $true, 1, $false, "4+27=", 2) | Out-Null
$d.Content.Find.Execute("79-36=", $true, $true, $false, $false, $false, $true, 1, $false, "60-55=", 2) | Out-Null
$d.Content.Find.Execute("65+20=", $true, $true, $false, $false, $false, $true, 1, $false, "55+28=", 2) | Out-Null
$d.Content.Find.Execute("28+53=", $true, $true, $false, $false, $false, $true, 1, $false, "14+60=", 2) | Out-Null
$d.Content.Find.Execute("36+6=", $true, $true, $false, $false, $false, $true, 1, $false, "50-5=", 2) | Out-Null
$d.Content.Find.Execute("79-56=", $true, $true, $false, $false, $false, $true, 1, $false, "96-58=", 2) | Out-Null
$d.Content.Find.Execute("32+3=", $true, $true, $false, $false, $false, $true, 1, $false, "69-40=", 2) | Out-Null
$d.Content.Find.Execute("15+14=", $true, $true, $false, $false, $false, $true, 1, $false, "98-43=", 2) | Out-Null
$d.Content.Find.Execute("36-15=", $true, $true, $false, $false, $false, $true, 1, $false, "68-64=", 2) | Out-Null
$d.Content.Find.Execute("21+57=", $true, $true, $false, $false, $false, $true, 1, $false, "61+27=", 2) | Out-Null
$d.Content.Find.Execute("21-20=", $true, $true, $false, $false, $false, $true, 1, $false, "64-18=", 2) | Out-Null
$d.Content.Find.Execute("21-6=", $true, $true, $false, $false, $false, $true, 1, $false, "1+5=", 2) | Out-Null
$d.Content.Find.Execute("57+8=", $true, $true, $false, $false, $false, $true, 1, $false, "63-27=", 2) | Out-Null
$d.Content.Find.Execute("13+48=", $true, $true, $false, $false, $false, $true, 1, $false, "63+5=", 2) | Out-Null
$d.Content.Find.Execute("97-77=", $true, $true, $false, $false, $false, $true, 1, $false, "49-46=", 2) | Out-Null
$d.Content.Find.Execute("59+8=", $true, $true, $false, $false, $false, $true, 1, $false, "95-34=", 2) | Out-Null
$d.Content.Find.Execute("92-32=", $true, $true, $false, $false, $false, $true, 1, $false, "17+3=", 2) | Out-Null
$d.Content.Find.Execute("27+11=", $true, $true, $false, $false, $false, $true, 1, $false, "34+62=", 2) | Out-Null
$d.Content.Find.Execute("18-1=", $true, $true, $false, $false, $false, $true, 1, $false, "19+24=", 2) | Out-Null
$d.Content.Find.Execute("43+3=", $true, $true, $false, $false, $false, $true, 1, $false, "49-18=", 2) | Out-Null
$d.Content.Find.Execute("97-0=", $true, $true, $false, $false, $false, $true, 1, $false, "20+65=", 2) | Out-Null
$d.Content.Find.Execute("82-61=", $true, $true, $false, $false, $false, $true, 1, $false, "69-59=", 2) | Out-Null
$d.Content.Find.Execute("59-36=", $true, $true, $false, $false, $false, $true, 1, $false, "44-15=", 2) | Out-Null
$d.Content.Find.Execute("29+14=", $true, $true, $false, $false, $false, $true, 1, $false, "16+35=", 2) | Out-Null
$d.Content.Find.Execute("12-5=", $true, $true, $false, $false, $false, $true, 1, $false, "74-11=", 2) | Out-Null
$d.Content.Find.Execute("31+35=", $true, $true, $false, $false, $false, $true, 1, $false, "5+4=", 2) | Out-Null
$d.Content.Find.Execute("69+15=", $true, $true, $false, $false, $false, $true, 1, $false, "9-8=", 2) | Out-Null
$d.Content.Find.Execute("49-14=", $true, $true, $false, $false, $false, $true, 1, $false, "88-78=", 2) | Out-Null
$d.Content.Find.Execute("13+6=", $true, $true, $false, $false, $false, $true, 1, $false, "7+36=", 2) | Out-Null
$d.Content.Find.Execute("35+41=", $true, $true, $false, $false, $false, $true, 1, $false, "29+24=", 2) | Out-Null
$d.Content.Find.Execute("31+20=", $true, $true, $false, $false, $false, $true, 1, $false, "20+27=", 2) | Out-Null
$d.Content.Find.Execute("13+41=", $true, $true, $false, $false, $false, $true, 1, $false, "11+10=", 2) | Out-Null
$d.Content.Find.Execute("23-0=", $true, $true, $false, $false, $false, $true, 1, $false, "43-20=", 2) | Out-Null
$d.Content.Find.Execute("33+16=", $true, $true, $false, $false, $false, $true, 1, $false, "58-38=", 2) | Out-Null
$d.Content.Find.Execute("65-59=", $true, $true, $false, $false, $false, $true, 1, $false, "63+13=", 2) | Out-Null
$d.Content.Find.Execute("12+22=", $true, $true, $false, $false, $false, $true, 1, $false, "27-20=", 2) | Out-Null
$d.Content.Find.Execute("65+6=", $true, $true, $false, $false, $false, $true, 1, $false, "75-25=", 2) | Out-Null
$d.Content.Find.Execute("20-13=", $true, $true, $false, $false, $false, $true, 1, $false, "14+4=", 2) | Out-Null
$d.Content.Find.Execute("90+7=", $true, $true, $false, $false, $false, $true, 1, $false, "24+17=", 2) | Out-Null
$d.Content.Find.Execute("61+5=", $true, $true, $false, $false, $false, $true, 1, $false, "30+7=", 2) | Out-Null
$d.Content.Find.Execute("0+15=", $true, $true, $false, $false, $false, $true, 1, $false, "20+35=", 2) | Out-Null
$d.Content.Find.Execute("1+51=", $true, $true, $false, $false, $false, $true, 1, $false, "94-83=", 2) | Out-Null
$d.Content.Find.Execute("99-50=", $true, $true, $false, $false, $false, $true, 1, $false, "4+40=", 2) | Out-Null
$d.Content.Find.Execute("20+48=", $true, $true, $false, $false, $false, $true, 1, $false, "91-3=", 2) | Out-Null
$d.Content.Find.Execute("97-97=", $true, $true, $false, $false, $false, $true, 1, $false, "16+75=", 2) | Out-Null
$d.Content.Find.Execute("40+21=", $true, $true, $false, $false, $false, $true, 1, $false, "51-19=", 2) | Out-Null
$d.Content.Find.Execute("37+50=", $true, $true, $false, $false, $false, $true, 1, $false, "0+13=", 2) | Out-Null
$d.Content.Find.Execute("40+25=", $true, $true, $false, $false, $false, $true, 1, $false, "30+56=", 2) | Out-Null
$d.Content.Find.Execute("41+58=", $true, $true, $false, $false, $false, $true, 1, $false, "20+59=", 2) | Out-Null
$d.Content.Find.Execute("61-6=", $true, $true, $false, $false, $false, $true, 1, $false, "63+36=", 2) | Out-Null
$d.Content.Find.Execute("92-12=", $true, $true, $false, $false, $false, $true, 1, $false, "35+32=", 2) | Out-Null
$d.Content.Find.Execute("70+19=", $true, $true, $false, $false, $false, $true, 1, $false, "53-19=", 2) | Out-Null
$d.Content.Find.Execute("40+35=", $true, $true, $false, $false, $false, $true, 1, $false, "37-29=", 2) | Out-Null
$d.Content.Find.Execute("77-38=", $true, $true, $false, $false, $false, $true, 1, $false, "1+49=", 2) | Out-Null
$d.Content.Find.Execute("33+38=", $true, $true, $false, $false, $false, $true, 1, $false, "63+1=", 2) | Out-Null
$d.Content.Find.Execute("64-35=", $true, $true, $false, $false, $false, $true, 1, $false, "11+42=", 2) | Out-Null
$d.Content.Find.Execute("57+29=", $true, $true, $false, $false, $false, $true, 1, $false, "95-85=", 2) | Out-Null
$d.Content.Find.Execute("84-81=", $true, $true, $false, $false, $false, $true, 1, $false, "21-17=", 2) | Out-Null
$d.Content.Find.Execute("99-4=", $true, $true, $false, $false, $false, $true, 1, $false, "26+1=", 2) | Out-Null
$d.Content.Find.Execute("35+25=", $true, $true, $false, $false, $false, $true, 1, $false, "34+59=", 2) | Out-Null
$d.Content.Find.Execute("36-33=", $true, $true, $false, $false, $false, $true, 1, $false, "69-27=", 2) | Out-Null
$d.Content.Find.Execute("44+53=", $true, $true, $false, $false, $false, $true, 1, $false, "20+22=", 2) | Out-Null
